$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '30.444.80'
Set-TextCell $ws.Range('E2') '  -0.44%  '
Set-TextCell $ws.Range('D3') '1.926.34'
Set-TextCell $ws.Range('E3') '  +3.86%  '
Set-TextCell $ws.Range('D4') '1.002'
Set-TextCell $ws.Range('E4') '  +0.28%  '
Set-TextCell $ws.Range('D5') '240.03'
Set-TextCell $ws.Range('E5') '  +2.61%  '
Set-TextCell $ws.Range('D6') '1.003'
Set-TextCell $ws.Range('E6') '  +0.34%  '
Set-TextCell $ws.Range('E7') '  +0.24%  '
Set-TextCell $ws.Range('B8') 'Cardano'
Set-TextCell $ws.Range('C8') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws.Range('D8') '0.2846'
Set-TextCell $ws.Range('E8') '  +3.48%  '
Set-TextCell $ws.Range('B9') 'Dogecoin'
Set-TextCell $ws.Range('C9') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws.Range('D9') '0.06562'
Set-TextCell $ws.Range('E9') '  +3.74%  '
Set-TextCell $ws.Range('B10') 'Solana'
Set-TextCell $ws.Range('C10') 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws.Range('D10') '19.22'
Set-TextCell $ws.Range('E10') '  +8.42%  '
Set-TextCell $ws.Range('B11') 'Litecoin'
Set-TextCell $ws.Range('C11') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws.Range('D11') '105.39'
Set-TextCell $ws.Range('E11') '  +24.62%  '
Set-TextCell $ws.Range('B12') 'WrappedEther'
Set-TextCell $ws.Range('C12') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range('D12') '1.935.40'
Set-TextCell $ws.Range('E12') '  +2.89%  '
Set-TextCell $ws.Range('B13') 'TRON'
Set-TextCell $ws.Range('C13') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws.Range('D13') '0.07586'
Set-TextCell $ws.Range('E13') '  +1.76%  '
Set-TextCell $ws.Range('B14') 'Polkadot'
Set-TextCell $ws.Range('C14') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range('D14') '5.121'
Set-TextCell $ws.Range('E14') '  +2.39%  '
Set-TextCell $ws.Range('B15') 'Polygon'
Set-TextCell $ws.Range('C15') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws.Range('D15') '0.6508'
Set-TextCell $ws.Range('E15') '  +3.93%  '
Set-TextCell $ws.Range('B16') 'BitcoinCash'
Set-TextCell $ws.Range('C16') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Range('D16') '298.76'
Set-TextCell $ws.Range('E16') '  +22.25%  '
Set-TextCell $ws.Range('B17') 'WrappedBTC'
Set-TextCell $ws.Range('C17') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws.Range('D17') '30.455.87'
Set-TextCell $ws.Range('E17') '  -0.26%  '
Set-TextCell $ws.Range('B18') 'Dai'
Set-TextCell $ws.Range('C18') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range('D18') '1.002'
Set-TextCell $ws.Range('E18') '  +0.19%  '
Set-TextCell $ws.Range('B19') 'Avalanche'
Set-TextCell $ws.Range('C19') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws.Range('D19') '12.96'
Set-TextCell $ws.Range('E19') '  +2.03%  '
Set-TextCell $ws.Range('B20') 'WrappedliquidstakedEther2.0'
Set-TextCell $ws.Range('C20') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws.Range('D20') '2.165.05'
Set-TextCell $ws.Range('E20') '  +3.88%  '
Set-TextCell $ws.Range('B21') 'ShibaInu'
Set-TextCell $ws.Range('C21') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws.Range('D21') '0.000007492'
Set-TextCell $ws.Range('E21') '  +2.15%  '
Set-TextCell $ws.Range('B22') 'BinanceUSD'
Set-TextCell $ws.Range('C22') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws.Range('D22') '1.001'
Set-TextCell $ws.Range('E22') '  +0.21%  '
Set-TextCell $ws.Range('B23') 'Uniswap'
Set-TextCell $ws.Range('C23') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws.Range('D23') '5.195'
Set-TextCell $ws.Range('E23') '  +5.15%  '
Set-TextCell $ws.Range('B24') 'Chainlink'
Set-TextCell $ws.Range('C24') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws.Range('D24') '6.299'
Set-TextCell $ws.Range('E24') '  +6.07%  '
Set-TextCell $ws.Range('B25') 'Cosmos'
Set-TextCell $ws.Range('C25') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws.Range('D25') '9.285'
Set-TextCell $ws.Range('E25') '  +1.50%  '
Set-TextCell $ws.Range('B26') 'Monero'
Set-TextCell $ws.Range('C26') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws.Range('D26') '165.28'
Set-TextCell $ws.Range('E26') '  +1.44%  '
Set-TextCell $ws.Range('B27') 'EthereumClassic'
Set-TextCell $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws.Range('D27') '19.63'
Set-TextCell $ws.Range('E27') '  +9.02%  '
Set-TextCell $ws.Range('B28') 'LidoDAOToken'
Set-TextCell $ws.Range('C28') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws.Range('D28') '2.032'
Set-TextCell $ws.Range('E28') '  +8.11%  '
Set-TextCell $ws.Range('B29') 'Stellar'
Set-TextCell $ws.Range('C29') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range('D29') '0.1124'
Set-TextCell $ws.Range('E29') '  +10.37%  '
Set-TextCell $ws.Range('B30') 'Toncoin'
Set-TextCell $ws.Range('C30') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range('D30') '1.349'
Set-TextCell $ws.Range('E30') '  -0.74%  '
Set-TextCell $ws.Range('B31') 'InternetComputer(DFINITY)'
Set-TextCell $ws.Range('C31') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range('D31') '4.114'
Set-TextCell $ws.Range('E31') '  +2.55%  '
Set-TextCell $ws.Range('B32') 'Filecoin'
Set-TextCell $ws.Range('C32') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range('D32') '3.922'
Set-TextCell $ws.Range('E32') '  +2.33%  '
Set-TextCell $ws.Range('B33') 'Hedera'
Set-TextCell $ws.Range('C33') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range('D33') '0.05019'
Set-TextCell $ws.Range('E33') '  +3.65%  '
Set-TextCell $ws.Range('B34') 'ImmutableX'
Set-TextCell $ws.Range('C34') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range('D34') '0.7367'
Set-TextCell $ws.Range('E34') '  +4.79%  '
Set-TextCell $ws.Range('B35') 'ARBITRUM'
Set-TextCell $ws.Range('C35') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range('D35') '1.145'
Set-TextCell $ws.Range('E35') '  +0.83%  '
Set-TextCell $ws.Range('B36') 'Frax'
Set-TextCell $ws.Range('C36') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell $ws.Range('D36') '1.001'
Set-TextCell $ws.Range('E36') '  +0.24%  '
Set-TextCell $ws.Range('B37') 'HuobiToken'
Set-TextCell $ws.Range('C37') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws.Range('D37') '2.721'
Set-TextCell $ws.Range('E37') '  +1.10%  '
Set-TextCell $ws.Range('B38') 'VeChain'
Set-TextCell $ws.Range('C38') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range('D38') '0.01960'
Set-TextCell $ws.Range('E38') '  +3.35%  '
Set-TextCell $ws.Range('B39') 'MXToken'
Set-TextCell $ws.Range('C39') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Range('D39') '2.699'
Set-TextCell $ws.Range('E39') '  +0.83%  '
Set-TextCell $ws.Range('B40') 'RenderToken'
Set-TextCell $ws.Range('C40') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range('D40') '2.017'
Set-TextCell $ws.Range('E40') '  +0.92%  '
Set-TextCell $ws.Range('B41') 'TrustWalletToken'
Set-TextCell $ws.Range('C41') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range('D41') '0.8738'
Set-TextCell $ws.Range('E41') '  -0.12%  '
Set-TextCell $ws.Range('B42') 'Quant'
Set-TextCell $ws.Range('C42') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws.Range('D42') '107.10'
Set-TextCell $ws.Range('E42') '  +0.33%  '
Set-TextCell $ws.Range('B43') 'FraxShare'
Set-TextCell $ws.Range('C43') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range('D43') '5.862'
Set-TextCell $ws.Range('E43') '  +5.71%  '
Set-TextCell $ws.Range('B44') 'PaxDollar'
Set-TextCell $ws.Range('C44') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws.Range('D44') '1.002'
Set-TextCell $ws.Range('E44') '  +0.27%  '
Set-TextCell $ws.Range('B45') 'Aave'
Set-TextCell $ws.Range('C45') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range('D45') '68.81'
Set-TextCell $ws.Range('E45') '  +9.66%  '
Set-TextCell $ws.Range('B46') 'TheSandbox'
Set-TextCell $ws.Range('C46') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws.Range('D46') '0.4135'
Set-TextCell $ws.Range('E46') '  +1.81%  '
Set-TextCell $ws.Range('B47') 'Aptos'
Set-TextCell $ws.Range('C47') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range('D47') '7.233'
Set-TextCell $ws.Range('E47') '  +0.55%  '
Set-TextCell $ws.Range('B48') 'EnergySwap'
Set-TextCell $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range('D48') '9.209'
Set-TextCell $ws.Range('E48') '  +7.54%  '
Set-TextCell $ws.Range('B49') 'Algorand'
Set-TextCell $ws.Range('C49') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range('D49') '0.1209'
Set-TextCell $ws.Range('E49') '  -0.21%  '
Set-TextCell $ws.Range('B50') 'Elrond'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell $ws.Range('D50') '34.62'
Set-TextCell $ws.Range('E50') '  +3.13%  '
Set-TextCell $ws.Range('B51') 'Cronos'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D51') '0.05636'
Set-TextCell $ws.Range('E51') '  +1.78%  '
